# Generate Report for handback
#
# For both locale sheets (zh-cn, de-de):
#   - Status cell (B2) flips from "Ready for handoff" to "Handed back: in sync with en-US"
#   - Two new cells appear: E2 (Latest Target File) and F2 (Latest Handback File),
#     mirroring the source file name (A2) and the handoff xlf name (C2), each carrying
#     a new "handback" hyperlink.
#   - The "Latest Handback DateTime" cell (G2) is stamped with the real handback time.

$wb = $excel.ActiveWorkbook

$locales = @(
    @{
        Sheet   = "zh-cn"
        XlfName = "d55136e5-0dee-46e6-9c1c-276ab72b7ddd.a48454d8b7dbdabf9ba764769a611f6a45bd1f93.zh-cn.xlf"
        MdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/c9e8cd4ad213164328b285c92494c0462d508866/e2e/d55136e5-0dee-46e6-9c1c-276ab72b7ddd.md"
        XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/97db59a6218f0e0fe4cd301e76270447871e9e24/ol-handback/OpenLocalizationTest/oltest.zh-cn/xinjiang/d55136e5-0dee-46e6-9c1c-276ab72b7ddd.a48454d8b7dbdabf9ba764769a611f6a45bd1f93.zh-cn.xlf"
        Handback = "2016-01-17 14:58:26"
    },
    @{
        Sheet   = "de-de"
        XlfName = "d55136e5-0dee-46e6-9c1c-276ab72b7ddd.a48454d8b7dbdabf9ba764769a611f6a45bd1f93.de-de.xlf"
        MdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/c9e8cd4ad213164328b285c92494c0462d508866/e2e/d55136e5-0dee-46e6-9c1c-276ab72b7ddd.md"
        XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d324f10569b9d5d40fbb61a6ca212ec8cbe6e59b/ol-handback/OpenLocalizationTest/oltest.de-de/xinjiang/d55136e5-0dee-46e6-9c1c-276ab72b7ddd.a48454d8b7dbdabf9ba764769a611f6a45bd1f93.de-de.xlf"
        Handback = "2016-01-17 14:58:44"
    }
)

$mdName = "d55136e5-0dee-46e6-9c1c-276ab72b7ddd.md"

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Status -> handed back
    $ws.Range("B2").Value = "Handed back: in sync with en-US"

    # Latest Target File (E2) - same file that was handed off, now target of handback
    $ws.Hyperlinks.Add($ws.Range("E2"), $loc.MdUrl, "", "", $mdName)
    $ws.Range("E2").Font.Underline = 2
    $ws.Range("E2").Font.Color = 15570276

    # Latest Handback File (F2) - the localized xlf coming back
    $ws.Hyperlinks.Add($ws.Range("F2"), $loc.XlfUrl, "", "", $loc.XlfName)
    $ws.Range("F2").Font.Underline = 2
    $ws.Range("F2").Font.Color = 15570276

    # Latest Handback DateTime (G2)
    $ws.Range("G2").Value = $loc.Handback
}
